# Apply the "Actualización automática de tasas-transfi.xlsx" update.

$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-of-the-day note (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nuevoTexto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.25 = 4687.5 pesos`n✅ 4687.5 pesos = 1.25 = 872.67 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $nuevoTexto

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 800
$wsTasas.Range("O10").Value = 3750
$wsTasas.Range("N12").Value = 3760
$wsTasas.Range("O12").Value = 700
